$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("Input"): restructure columns from the old 17-column layout
# (A:Q) to the new 16-column standard-template layout (A:P).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Wipe all existing content/formatting so we can rebuild from scratch.
$ws1.Cells.Clear()

# --- Header row (row 1) - plain, unstyled text -----------------------------
# (Cells are set individually/scalar - assigning an array to a Range.Value
# in one shot confuses this interpreter's output handling.)
$ws1.Range("A1").Value = "발주일자"
$ws1.Range("B1").Value = "납기일자"
$ws1.Range("C1").Value = "거래처명"
$ws1.Range("D1").Value = "거래처 이메일"
$ws1.Range("E1").Value = "납품처명"
$ws1.Range("F1").Value = "납품처 이메일"
$ws1.Range("G1").Value = "프로젝트명"
$ws1.Range("H1").Value = "대분류"
$ws1.Range("I1").Value = "중분류"
$ws1.Range("J1").Value = "소분류"
$ws1.Range("K1").Value = "품목명"
$ws1.Range("L1").Value = "규격"
$ws1.Range("M1").Value = "수량"
$ws1.Range("N1").Value = "단가"
$ws1.Range("O1").Value = "총금액"
$ws1.Range("P1").Value = "비고"

# --- Data rows ---------------------------------------------------------------
# Columns: A 발주일자 | B 납기일자 | C 거래처명 | D 거래처 이메일 | E 납품처명
#          F 납품처 이메일 | G 프로젝트명 | H 대분류 | I 중분류 | J 소분류
#          K 품목명 | L 규격 | M 수량 | N 단가 | O 총금액 | P 비고

$rows = @(
    @{ Row=2; OrderDate="2025-08-30"; DueDate="2025-09-27"; Vendor="유니모터스"; VendorEmail="유니모터스@example.com"; Site="힐스테이트 도곡동1차"; SiteEmail="delivery@example.com"; Project="힐스테이트 도곡동1차"; Major="4. 장비비"; Middle="1) 장비비"; Minor="기타"; Item="자재하차"; Spec="KS규격-1"; Qty=2; UnitPrice=80000; Total=176000 },
    @{ Row=3; OrderDate="2025-08-25"; DueDate="2025-09-30"; Vendor="유니모터스"; VendorEmail="유니모터스@example.com"; Site="힐스테이트 도곡동1차"; SiteEmail="delivery@example.com"; Project="힐스테이트 도곡동1차"; Major="4. 장비비"; Middle="1) 장비비"; Minor="기타"; Item="렌탈이동"; Spec="KS규격-2"; Qty=3; UnitPrice=80000; Total=264000 },
    @{ Row=4; OrderDate="2025-09-13"; DueDate="2025-09-03"; Vendor="유니모터스"; VendorEmail="유니모터스@example.com"; Site="힐스테이트 도곡동1차"; SiteEmail="delivery@example.com"; Project="힐스테이트 도곡동1차"; Major="4. 장비비"; Middle="1) 장비비"; Minor="기타"; Item="자재하차"; Spec="KS규격-3"; Qty=1; UnitPrice=80000; Total=88000 },
    @{ Row=5; OrderDate="2025-09-03"; DueDate="2025-10-16"; Vendor="유니모터스"; VendorEmail="유니모터스@example.com"; Site="힐스테이트 도곡동1차"; SiteEmail="delivery@example.com"; Project="힐스테이트 도곡동1차"; Major="4. 장비비"; Middle="1) 장비비"; Minor="기타"; Item="자재하차"; Spec="KS규격-4"; Qty=1; UnitPrice=80000; Total=88000 },
    @{ Row=6; OrderDate="2025-08-25"; DueDate="2025-10-01"; Vendor="유니모터스"; VendorEmail="유니모터스@example.com"; Site="힐스테이트 도곡동1차"; SiteEmail="delivery@example.com"; Project="힐스테이트 도곡동1차"; Major="4. 장비비"; Middle="1) 장비비"; Minor="기타"; Item="자재하차"; Spec="KS규격-5"; Qty=1; UnitPrice=80000; Total=88000 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Date-like text must be forced to Text format first, otherwise Excel's
    # smart-parsing will silently convert "2025-08-30" into a date serial.
    $dateRange = $ws1.Range("A" + $rowNum + ":B" + $rowNum)
    $dateRange.NumberFormat = "@"
    $ws1.Range("A" + $rowNum).Value = $r.OrderDate
    $ws1.Range("B" + $rowNum).Value = $r.DueDate
    $dateRange.Style = "Normal"

    $ws1.Range("C" + $rowNum).Value = $r.Vendor
    $ws1.Range("D" + $rowNum).Value = $r.VendorEmail
    $ws1.Range("E" + $rowNum).Value = $r.Site
    $ws1.Range("F" + $rowNum).Value = $r.SiteEmail
    $ws1.Range("G" + $rowNum).Value = $r.Project
    $ws1.Range("H" + $rowNum).Value = $r.Major
    $ws1.Range("I" + $rowNum).Value = $r.Middle
    $ws1.Range("J" + $rowNum).Value = $r.Minor
    $ws1.Range("K" + $rowNum).Value = $r.Item
    $ws1.Range("L" + $rowNum).Value = $r.Spec
    $ws1.Range("M" + $rowNum).Value = $r.Qty
    $ws1.Range("N" + $rowNum).Value = $r.UnitPrice
    $ws1.Range("O" + $rowNum).Value = $r.Total
    # Column P (비고) intentionally left blank for every data row.
}

# ---------------------------------------------------------------------------
# Sheet 2 ("갑지") and Sheet 3 ("을지"): the only change is that the
# previously-empty "비고" cells (I2:I6), which held an empty inline string,
# are now fully empty (no cell content at all).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("I2:I6").ClearContents()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("I2:I6").ClearContents()
